$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
    }
    return $found
}

# 1. Notebook / datasets sentence split
Replace-Text "respectively, that were later merged and further analyzed by district and school." "respectively. These datasets were then merged and further analyzed by district and school."

# 2. "of students met this requirement, respectively."
Replace-Text "of students met this requirement, respectively." "of the students met this requirement for math and reading, respectively."

# 3. "analysis also included average math and reading scores per school"
Replace-Text "analysis also included average math and reading scores per school" "analysis also included the average math and reading scores per school"

# 5. "Wilson and Pena high schools ... Johnson high schools."
Replace-Text "Wilson and Pena high schools, whereas the top 5 lowest performing schools include Rodriguez, Figueroa, Huang, Hernandez and Johnson high schools." "Wilson and Pena High Schools, whereas the top 5 lowest performing schools include Rodriguez, Figueroa, Huang, Hernandez and Johnson High Schools."

# 9. "Hernandez and Huang high schools had among the highest budget per student, they were among the lowest performing school by % Overall Passing."
Replace-Text "Hernandez and Huang high schools had among the highest budget per student, they were among the lowest performing school by % Overall Passing." "Hernandez and Huang High Schools had among the highest budget per student, they were among the lowest performing schools by % Overall Passing."

# 10. "Regardless of school, there seems to be a trend where the" (non-bold) -> "...seemed..."
Replace-Text "Regardless of school, there seems to be a trend where the" "Regardless of school, there seemed to be a trend where the"

# "higher the spending ranges per student, the lower is the % Overall Passing" (bold run) -> "...the lower was..." (do not cross into the trailing bold ". " run)
Replace-Text "higher the spending ranges per student, the lower is the % Overall Passing" "higher the spending ranges per student, the lower was the % Overall Passing"

# 11. "are " (bold run, "...math and reading scores are also associated...") -> "were "
Replace-Text "math and reading scores are also associated with lower spending ranges per student." "math and reading scores were also associated with lower spending ranges per student."

# 12. Insert new empty NoSpacing paragraph after "... A summary sample of the results is shown below." (second occurrence - after Wilson High School sentence)
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Wilson High School was the only large school with one of the highest % Overall Passing. A summary sample of the results is shown below.")
if ($found) {
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1)
    $p = $rng.Paragraphs(1)
    $p.Style = "NoSpacing"
} else {
    Write-Host "NOT FOUND: Wilson High School summary sample paragraph"
}

# 13. "smaller the school size, the higher is the % Overall Passing" (bold run) -> "... the higher was the % Overall Passing" (keep bold boundary intact, do not bleed into the preceding/following non-bold text)
Replace-Text "smaller the school size, the higher is the % Overall Passing" "smaller the school size, the higher was the % Overall Passing"

# 14. "Overall, data suggest that Charter schools perform petter than District schools in math, it is slightly better in reading, and the % Overall Passing is higher for Charter schools."
Replace-Text "Overall, data suggest that Charter schools perform petter than District schools in math, it is slightly better in reading, and the % Overall Passing is higher for Charter schools." "Overall, data suggested that Charter schools performed petter than District schools in math, and their % Overall Passing was also higher compared to District schools."

# 15. "Higher budget and spending values per student are not necessarily associated with better performance."
Replace-Text "budget and spending values per student are not necessarily associated with better performance." "budget and spending ranges per student were not necessarily associated with better performance."

Replace-Text "several of the high schools with the highest budget and spending values per student had lower % Overall Passing" "several of the high schools with the highest budget and spending ranges per student had lower % Overall Passing"

Replace-Text "The opposite was true for schools that had lower budgets and spending values per student." "The opposite was true for most schools that had lower budgets and spending values per student."
